$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.757.30"
$ws.Range("D3").Value = "2.630.53"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.574"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").Value = "2.637.08"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.97%  "
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.335"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "3.092.69"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "58.732.27"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000137"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.21%  "
$ws.Range("D18").Value = "2.634.52"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "348.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("E20").Value = "  -3.76%  "
$ws.Range("E21").Value = "  -3.28%  "
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.415"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.01%  "
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("D28").Value = "0.0₃0805"
$ws.Range("E28").Value = "  -4.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("E33").Value = "  -0.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.984"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.844"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.36%  "
$ws.Range("E41").Value = "  -3.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "280.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.997"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0986"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.598"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0525"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("D49").Value = "1.994.93"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("E50").Value = "  -2.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.77%  "
